$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34").Value = 43237
$ws.Range("B34").Value = "Suite refactoring et débugging"
$ws.Range("C34").Value = 12

$ws.Range("A35").Value = 43238
$ws.Range("B35").Value = "Suite refactoring et débugging"
$ws.Range("C35").Value = 1.5

$ws.Range("B36").Select()
